$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----- Header text updates (new crime data collection week) -----
# A8 holds "Volume 31   Number  30" -> bump the issue number to 31
$ws.Range("A8").Characters(21, 2).Text = "31"
# C9 holds "Report Covering the Week  7/22/2024  Through  7/28/2024"
# -> advance the reporting week by one (7/29/2024 through 8/4/2024)
$ws.Range("C9").Characters(27, 9).Text = "7/29/2024"
$ws.Range("C9").Characters(47, 9).Text = "8/4/2024"

# Number formats matching existing styles already present in this workbook
# (reusing them avoids creating duplicate style/numFmt entries)
$fmtInt  = "#,##0"
$fmtDec1 = '#,##0.0;"-"#,##0.0'

# Donor cells (row 14, untouched by this edit) already carrying the special
# placeholder text + style used for N/A-type cells ("0" and "***.*")
$donorZero = $ws.Range("C14")   # text "0"
$donorNA   = $ws.Range("E14")   # text "***.*"

# ----- Row 16 -----
$ws.Range("C16").Value = 1
$ws.Range("D16").NumberFormat = $fmtInt
$ws.Range("D16").Value = 1
$ws.Range("E16").NumberFormat = $fmtDec1
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 150
$ws.Range("I16").Value = 29
$ws.Range("J16").Value = 9
$ws.Range("K16").Value = 222.222222222222
$ws.Range("L16").Value = 222.222222222222
$ws.Range("M16").Value = 81.25
$ws.Range("N16").Value = -74.561403508771

# ----- Row 17 -----
$donorZero.Copy($ws.Range("C17"))
$ws.Range("F17").Value = 2
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 150
$ws.Range("N17").Value = -61.538461538461

# ----- Row 19 -----
$ws.Range("C19").NumberFormat = $fmtInt
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = -66.666666666666
$ws.Range("F19").Value = 2
$ws.Range("G19").Value = 11
$ws.Range("H19").Value = -81.818181818181
$ws.Range("I19").Value = 29
$ws.Range("J19").Value = 30
$ws.Range("K19").Value = -3.333333333333
$ws.Range("L19").Value = 81.25
$ws.Range("M19").Value = -35.555555555555
$ws.Range("N19").Value = -72.380952380952

# ----- Row 21 -----
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = 4
$ws.Range("E21").Value = -50
$ws.Range("G21").Value = 13
$ws.Range("H21").Value = -30.769230769230
$ws.Range("I21").Value = 70
$ws.Range("J21").Value = 47
$ws.Range("K21").Value = 48.936170212766
$ws.Range("L21").Value = 79.487179487179
$ws.Range("M21").Value = 1.449275362318
$ws.Range("N21").Value = -74.637681159420

# ----- Row 24 -----
$ws.Range("C24").NumberFormat = $fmtInt
$ws.Range("C24").Value = 1
$donorZero.Copy($ws.Range("D24"))
$donorNA.Copy($ws.Range("E24"))
$ws.Range("F24").Value = 2
$ws.Range("G24").Value = 9
$ws.Range("H24").Value = -77.777777777777
$ws.Range("I24").Value = 20
$ws.Range("K24").Value = -20
$ws.Range("L24").Value = 11.111111111111
$ws.Range("M24").Value = -56.521739130434

# ----- Row 26 -----
$donorZero.Copy($ws.Range("D26"))
$donorNA.Copy($ws.Range("E26"))
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = -60
$ws.Range("M26").Value = 64.285714285714

# ----- Row 28 -----
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 2
$ws.Range("I28").Value = 7
$ws.Range("J28").Value = 14
$ws.Range("K28").Value = -50
$ws.Range("L28").Value = 0
